$wb = $excel.ActiveWorkbook

# Locate the existing TestData sheet and duplicate it; Excel places the copy
# immediately before the source sheet, so it becomes the new first tab.
$testData = $wb.Worksheets.Item("TestData")
$testData.Copy($testData)

# NOTE: inserting a sheet shifts worksheet indices/handles, so re-resolve
# every sheet reference by name afterwards instead of reusing old COM
# references created before the Copy().
$keyTestData = $wb.Worksheets.Item("TestData (2)")
$keyTestData.Name = "KeyTestData"

# The new sheet holds the "keyed" variant of the test data: the first
# column's type header becomes "int;key" instead of plain "int" so the
# importer knows column A is the record key.
$keyTestData.Range("A1").Value = "int;key"

# Restore a sane per-sheet selection on TestData (it keeps its own cursor
# position, independent from the new sheet).
$testData = $wb.Worksheets.Item("TestData")
$testData.Activate()
[void]$testData.Range("D12").Select()

# Leave MapData untouched, then make KeyTestData the active/selected tab,
# matching the new workbook view after the edit.
$keyTestData = $wb.Worksheets.Item("KeyTestData")
$keyTestData.Activate()
[void]$keyTestData.Range("D21").Select()
